$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 15 (C15:H15) with 5s
$ws.Range("C15:H15").Value = 5

# Row 20 needs an extra I20 cell (matching the existing H20 formatting)
# before filling in the scores, since I20 was previously empty/unstyled.
$ws.Range("H20").Copy()
$ws.Range("I20").PasteSpecial(-4122)

# Fill in row 20 (C20:I20) with 5s
$ws.Range("C20:I20").Value = 5

# Update the selection to match the new active cell / selected range
$ws.Range("C20:I20").Select()
